$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'rühren',
    'münzen',
    'grüßen',
    'jagen',
    'schnellen',
    'tragen',
    'seufzen',
    'ändern',
    'lesen',
    'bluten',
    'gnaden',
    'zünden',
    'saufen',
    'werfen',
    'hindern',
    'nutzen',
    'warnen',
    'reizen',
    'plaudern',
    'kratzen',
    'enden',
    'hauen',
    'hassen',
    'malen',
    'lügen',
    'pflanzen',
    'dienen',
    'suchen',
    'danken',
    'äußern',
    'tanzen',
    'spüren',
    'flüchten',
    'löschen',
    'binden',
    'boxen',
    'fallen',
    'kehren',
    'kürzen',
    'ehren',
    'garen',
    'achten',
    'zielen',
    'hören',
    'graben',
    'irren',
    'stecken',
    'rücken',
    'fehlen',
    'wüten',
    'schleppen',
    'wenden',
    'stören',
    'klingen',
    'klagen',
    'treiben',
    'werden',
    'lockern',
    'geben',
    'siegen',
    'knarren',
    'zahlen',
    'sperren',
    'sparen',
    'helfen',
    'wirken',
    'feiern',
    'teilen',
    'mögen',
    'flehen',
    'fangen',
    'jubeln',
    'sagen',
    'folgen',
    'kämpfen',
    'regeln',
    'sorgen',
    'fällen',
    'husten',
    'schulden',
    'herzen',
    'heilen',
    'schlafen',
    'erben',
    'zeigen',
    'kichern',
    'liefern',
    'sichern',
    'biegen',
    'knurren',
    'planen',
    'kümmern',
    'altern',
    'bleiben',
    'mauern',
    'schenken',
    'trennen',
    'rechnen',
    'bitten',
    'klettern',
    'betteln',
    'quälen',
    'brauchen',
    'scheitern',
    'messen',
    'pissen',
    'schlucken',
    'lieben',
    'träumen',
    'stammen',
    'bellen',
    'brauen',
    'töten',
    'zeugen',
    'sterben',
    'feuern',
    'drehen',
    'merken',
    'loben',
    'schwächen',
    'spielen',
    'freuen',
    'rufen',
    'wundern',
    'rauchen',
    'wiegen',
    'filmen',
    'nähen',
    'pfeifen',
    'wehren',
    'decken',
    'führen',
    'dringen',
    'fahren',
    'weinen',
    'gelten',
    'platzen',
    'greifen',
    'läuten',
    'ärgern',
    'wehtun',
    'regnen',
    'spinnen',
    'liegen',
    'lösen',
    'heulen',
    'formen',
    'arten',
    'schmecken',
    'herrschen',
    'scheinen',
    'machen',
    'bilden',
    'tollen',
    'trauen',
    'strahlen',
    'wachsen',
    'weichen',
    'rauben',
    'kranken',
    'heben',
    'lächeln',
    'ächzen',
    'fischen',
    'runden',
    'schwören',
    'schreiten',
    'schwingen',
    'meinen',
    'backen',
    'schütteln',
    'fühlen',
    'zeichnen',
    'räumen',
    'sinken',
    'stimmen',
    'rasen',
    'schrecken',
    'bergen',
    'fließen',
    'bauen',
    'füttern',
    'gründen',
    'drohen',
    'sammeln',
    'stehlen',
    'tropfen',
    'stoßen',
    'sprengen',
    'kosten',
    'streichen',
    'zögern'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}
